$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Affiliate disclosure in introduction*") {
        $p.Range.Delete()
        break
    }
}
